$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = [double]"24.85000000000045"
$ws.Range("H2").Value = [double]"2.306956934286039e-16"
$ws.Range("K2").Value = [double]"44.19802032410474"
$ws.Range("L2").Value = "[32.43009022766118, 55.965950420548296]"
$ws.Range("M2").Value = [double]"1.608713162681852e-12"
$ws.Range("N2").Value = [double]"1.608713162681852e-12"
$ws.Range("O2").Value = [double]"2.018921405009657"
$ws.Range("P2").Value = "[1.7547634641672714, 2.283079345852042]"
$ws.Range("S2").Value = [double]"61.77431356931264"
$ws.Range("T2").Value = "[55.52944715977745, 68.01917997884783]"
$ws.Range("W2").Value = [double]"16.86516516516547"
$ws.Range("X2").Value = [double]"15.8204204204207"
$ws.Range("Y2").Value = [double]"17.90990990991023"

# Row 3 updates
$ws.Range("B3").Value = [double]"1"
$ws.Range("E3").Value = [double]"24.11000000000033"
$ws.Range("H3").Value = [double]"2.306956934286039e-16"
$ws.Range("K3").Value = [double]"48.44185049171178"
$ws.Range("L3").Value = "[36.552364359586775, 60.331336623836776]"
$ws.Range("M3").Value = [double]"4.662936703425657e-14"
$ws.Range("N3").Value = [double]"9.325873406851315e-14"
$ws.Range("O3").Value = [double]"1.150973885098963"
$ws.Range("P3").Value = "[0.8868159442565782, 1.4151318259413488]"
$ws.Range("Q3").Value = [double]"1.110223024625157e-15"
$ws.Range("R3").Value = [double]"1.110223024625157e-15"
$ws.Range("S3").Value = [double]"66.55196693381382"
$ws.Range("T3").Value = "[59.87940612818599, 73.22452773944164]"
$ws.Range("W3").Value = [double]"19.69345345345372"
$ws.Range("X3").Value = [double]"18.67981981982007"
$ws.Range("Y3").Value = [double]"20.70708708708737"
